# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    # Force the value to be stored as text (not auto-converted to a number),
    # then restore the default "Normal" style so no extra formatting sticks.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

function Set-Price($row, $value) {
    Set-Text $row 4 $value
}

function Set-Volume($row, $value) {
    Set-Text $row 5 $value
}

# Row 2 - Bitcoin
Set-Price 2 "71.737.39"
Set-Volume 2 "  +2.94%  "

# Row 3 - Ethereum
Set-Price 3 "3.631.57"
Set-Volume 3 "  +6.64%  "

# Row 4 - TetherUSD
Set-Volume 4 "  +0.01%  "

# Row 5 - BNB
Set-Price 5 "587.77"
Set-Volume 5 "  +0.08%  "

# Row 6 - Solana
Set-Price 6 "180.63"
Set-Volume 6 "  -0.38%  "

# Row 7 - LidoStakedEther
Set-Price 7 "3.618.08"
Set-Volume 7 "  +6.43%  "

# Row 8 - XRP
Set-Price 8 "0.615"
Set-Volume 8 "  +2.54%  "

# Row 9 - USDC
Set-Volume 9 "  +0.00%  "

# Row 10 - Dogecoin
Set-Price 10 "0.202"
Set-Volume 10 "  -0.75%  "

# Row 11 - Cardano
Set-Price 11 "0.607"
Set-Volume 11 "  +2.23%  "

# Row 12 - Avalanche
Set-Price 12 "49.72"
Set-Volume 12 "  +2.34%  "

# Row 13 - ShibaInu
Set-Price 13 "0.0000286"
Set-Volume 13 "  -0.73%  "

# Row 14 - BitcoinCash
Set-Price 14 "683.18"
Set-Volume 14 "  -0.25%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-Price 15 "4.220.10"
Set-Volume 15 "  +6.52%  "

# Row 16 - Polkadot
Set-Price 16 "9.02"
Set-Volume 16 "  +3.84%  "

# Row 17 - WrappedEther
Set-Price 17 "3.655.46"
Set-Volume 17 "  +7.01%  "

# Row 18 - WrappedBTC
Set-Price 18 "71.858.63"
Set-Volume 18 "  +3.00%  "

# Row 19 - TRON
Set-Volume 19 "  +1.60%  "

# Row 20 - Chainlink
Set-Price 20 "18.29"
Set-Volume 20 "  +2.81%  "

# Row 21 - Uniswap
Set-Price 21 "11.61"
Set-Volume 21 "  +2.07%  "

# Row 22 - Polygon
Set-Price 22 "0.938"
Set-Volume 22 "  +2.58%  "

# Row 23 - Toncoin
Set-Price 23 "5.90"
Set-Volume 23 "  +10.01%  "

# Row 24 - InternetComputer(DFINITY)
Set-Price 24 "17.75"
Set-Volume 24 "  +2.30%  "

# Row 25 - Litecoin
Set-Price 25 "103.29"
Set-Volume 25 "  +0.12%  "

# Row 26 - PancakeSwap
Set-Price 26 "4.01"
Set-Volume 26 "  +1.37%  "

# Row 27 - ImmutableX
Set-Price 27 "2.85"
Set-Volume 27 "  +4.58%  "

# Row 28 - RenderToken
Set-Price 28 "9.98"
Set-Volume 28 "  +1.99%  "

# Row 29 - EthereumClassic
Set-Price 29 "35.10"
Set-Volume 29 "  +3.01%  "

# Row 30 - Filecoin
Set-Price 30 "9.22"
Set-Volume 30 "  +4.21%  "

# Row 31 - NEARProtocol
Set-Price 31 "7.34"
Set-Volume 31 "  +5.01%  "

# Row 32 - dogwifhat
Set-Volume 32 "  +15.23%  "

# Row 33 - Bittensor
Set-Price 33 "587.65"
Set-Volume 33 "  +5.54%  "

# Row 34 - Cosmos
Set-Price 34 "11.32"
Set-Volume 34 "  +1.41%  "

# Row 36 - OKB
Set-Volume 36 "  +1.34%  "

# Row 37 - Dai
Set-Volume 37 "  +0.01%  "

# Row 38 - Maker
Set-Price 38 "3.683.39"
Set-Volume 38 "  +0.53%  "

# Row 39 - Kaspa
Set-Price 39 "0.142"
Set-Volume 39 "  +0.77%  "

# Row 40 - InjectiveProtocol
Set-Price 40 "35.64"
Set-Volume 40 "  -0.79%  "

# Row 41 - PEPE (value contains subscript-3 unicode char U+2083)
$sub3 = [char]0x2083
$pepePrice = "0.0{0}0763" -f $sub3
Set-Price 41 $pepePrice
Set-Volume 41 "  +3.72%  "

# Rows 42/43 swap: Stacks moves to row 42, VeChain moves to row 43
Set-Text 42 2 "Stacks"
Set-Text 42 3 "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-Price 42 "3.41"
Set-Volume 42 "  +3.45%  "

Set-Text 43 2 "VeChain"
Set-Text 43 3 "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-Price 43 "0.0469"
Set-Volume 43 "  +8.15%  "

# Row 44 - Fetch.AI
Set-Price 44 "2.78"
Set-Volume 44 "  +2.64%  "

# Row 45 - TheGraph
Set-Price 45 "0.346"
Set-Volume 45 "  +1.53%  "

# Row 46 - ApeXProtocol
Set-Price 46 "3.42"
Set-Volume 46 "  +1.73%  "

# Row 47 - ThetaToken
Set-Price 47 "2.80"
Set-Volume 47 "  +4.41%  "

# Row 48 - Stellar
Set-Volume 48 "  +2.80%  "

# Row 49 - Mantle
Set-Volume 49 "  +3.12%  "

# Row 50 - FirstDigitalUSD
Set-Volume 50 "  -0.02%  "

# Row 51 - Monero
Set-Price 51 "131.38"
Set-Volume 51 "  +1.11%  "
